$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("September")

$ws.Range("B2").Value = 1586
$ws.Range("C2").Value = 1178
$ws.Range("D2").Value = 408
$ws.Range("E2").Value = "We borrowerd more than we lent"
$ws.Range("G2").Value = "1.35 : 1"

$ws.Range("B3").Value = 820
$ws.Range("C3").Value = 537
$ws.Range("D3").Value = 283
$ws.Range("E3").Value = "We borrowerd more than we lent"
$ws.Range("G3").Value = "1.53 : 1"

$ws.Range("B4").Value = 1263
$ws.Range("C4").Value = 1366
$ws.Range("D4").Value = -103
$ws.Range("F4").Value = "We lent more than we borrowed"
$ws.Range("G4").Value = "0.92 : 1"

$ws.Range("B5").Value = 124
$ws.Range("C5").Value = 147
$ws.Range("D5").Value = -23
$ws.Range("F5").Value = "We lent more than we borrowed"
$ws.Range("G5").Value = "0.84 : 1"

$ws.Range("B6").Value = 1315
$ws.Range("C6").Value = 1519
$ws.Range("D6").Value = -204
$ws.Range("F6").Value = "We lent more than we borrowed"
$ws.Range("G6").Value = "0.87 : 1"

$ws.Range("B7").Value = 172
$ws.Range("C7").Value = 225
$ws.Range("D7").Value = -53
$ws.Range("F7").Value = "We lent more than we borrowed"
$ws.Range("G7").Value = "0.76 : 1"

$ws.Range("B8").Value = 163
$ws.Range("C8").Value = 140
$ws.Range("D8").Value = 23
$ws.Range("E8").Value = "We borrowerd more than we lent"
$ws.Range("G8").Value = "1.16 : 1"

$ws.Range("B9").Value = 37
$ws.Range("C9").Value = 67
$ws.Range("D9").Value = -30
$ws.Range("F9").Value = "We lent more than we borrowed"
$ws.Range("G9").Value = "0.55 : 1"

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 58
$ws.Range("D10").Value = -58
$ws.Range("F10").Value = "We lent more than we borrowed"
$ws.Range("G10").Value = "0.00 : 1"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 37
$ws.Range("D12").Value = -26
$ws.Range("F12").Value = "We lent more than we borrowed"
$ws.Range("G12").Value = "0.30 : 1"

$ws.Range("B13").Value = 148
$ws.Range("C13").Value = 111
$ws.Range("D13").Value = 37
$ws.Range("E13").Value = "We borrowerd more than we lent"
$ws.Range("G13").Value = "1.33 : 1"

$ws.Range("B14").Value = 108
$ws.Range("C14").Value = 297
$ws.Range("D14").Value = -189
$ws.Range("F14").Value = "We lent more than we borrowed"
$ws.Range("G14").Value = "0.36 : 1"

$ws.Range("B15").Value = 127
$ws.Range("C15").Value = 97
$ws.Range("D15").Value = 30
$ws.Range("E15").Value = "We borrowerd more than we lent"
$ws.Range("G15").Value = "1.31 : 1"

$ws.Range("B16").Value = 35
$ws.Range("C16").Value = 158
$ws.Range("D16").Value = -123
$ws.Range("F16").Value = "We lent more than we borrowed"
$ws.Range("G16").Value = "0.22 : 1"

$ws.Range("B17").Value = 741
$ws.Range("C17").Value = 517
$ws.Range("D17").Value = 224
$ws.Range("E17").Value = "We borrowerd more than we lent"
$ws.Range("G17").Value = "1.43 : 1"

$ws.Range("B18").Value = 88
$ws.Range("C18").Value = 107
$ws.Range("D18").Value = -19
$ws.Range("F18").Value = "We lent more than we borrowed"
$ws.Range("G18").Value = "0.82 : 1"

$ws.Range("B19").Value = 629
$ws.Range("C19").Value = 499
$ws.Range("D19").Value = 130
$ws.Range("E19").Value = "We borrowerd more than we lent"
$ws.Range("G19").Value = "1.26 : 1"

$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 74
$ws.Range("D20").Value = -70
$ws.Range("F20").Value = "We lent more than we borrowed"
$ws.Range("G20").Value = "0.05 : 1"

$ws.Range("B21").Value = 436
$ws.Range("C21").Value = 433
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = "We borrowerd more than we lent"
$ws.Range("G21").Value = "1.01 : 1"

$ws.Range("B22").Value = 38
$ws.Range("C22").Value = 90
$ws.Range("D22").Value = -52
$ws.Range("F22").Value = "We lent more than we borrowed"
$ws.Range("G22").Value = "0.42 : 1"

$ws.Range("B23").Value = 724
$ws.Range("C23").Value = 388
$ws.Range("D23").Value = 336
$ws.Range("E23").Value = "We borrowerd more than we lent"
$ws.Range("G23").Value = "1.87 : 1"

$ws.Range("B24").Value = 1997
$ws.Range("C24").Value = 1377
$ws.Range("D24").Value = 620
$ws.Range("E24").Value = "We borrowerd more than we lent"
$ws.Range("G24").Value = "1.45 : 1"

$ws.Range("B25").Value = 176
$ws.Range("C25").Value = 368
$ws.Range("D25").Value = -192
$ws.Range("F25").Value = "We lent more than we borrowed"
$ws.Range("G25").Value = "0.48 : 1"

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = 176
$ws.Range("C27").Value = 215
$ws.Range("D27").Value = -39
$ws.Range("F27").Value = "We lent more than we borrowed"
$ws.Range("G27").Value = "0.82 : 1"

$ws.Range("B28").Value = 51
$ws.Range("C28").Value = 72
$ws.Range("D28").Value = -21
$ws.Range("F28").Value = "We lent more than we borrowed"
$ws.Range("G28").Value = "0.71 : 1"

$ws.Range("B29").Value = 563
$ws.Range("C29").Value = 482
$ws.Range("D29").Value = 81
$ws.Range("E29").Value = "We borrowerd more than we lent"
$ws.Range("G29").Value = "1.17 : 1"

$ws.Range("B30").Value = 39
$ws.Range("C30").Value = 39
$ws.Range("D30").Value = 0
$ws.Range("G30").Value = "1.00 : 1"

$ws.Range("B31").Value = 76
$ws.Range("C31").Value = 259
$ws.Range("D31").Value = -183
$ws.Range("F31").Value = "We lent more than we borrowed"
$ws.Range("G31").Value = "0.29 : 1"

$ws.Range("B32").Value = 457
$ws.Range("C32").Value = 546
$ws.Range("D32").Value = -89
$ws.Range("F32").Value = "We lent more than we borrowed"
$ws.Range("G32").Value = "0.84 : 1"

$ws.Range("B33").Value = 379
$ws.Range("C33").Value = 492
$ws.Range("D33").Value = -113
$ws.Range("F33").Value = "We lent more than we borrowed"
$ws.Range("G33").Value = "0.77 : 1"

$ws.Range("B34").Value = 177
$ws.Range("C34").Value = 124
$ws.Range("D34").Value = 53
$ws.Range("E34").Value = "We borrowerd more than we lent"
$ws.Range("G34").Value = "1.43 : 1"

$ws.Range("B35").Value = 1060
$ws.Range("C35").Value = 1016
$ws.Range("D35").Value = 44
$ws.Range("E35").Value = "We borrowerd more than we lent"
$ws.Range("G35").Value = "1.04 : 1"

$ws.Range("B36").Value = 179
$ws.Range("C36").Value = 469
$ws.Range("D36").Value = -290
$ws.Range("F36").Value = "We lent more than we borrowed"
$ws.Range("G36").Value = "0.38 : 1"

$ws.Range("B37").Value = 549
$ws.Range("C37").Value = 312
$ws.Range("D37").Value = 237
$ws.Range("E37").Value = "We borrowerd more than we lent"
$ws.Range("G37").Value = "1.76 : 1"

$ws.Range("B38").Value = 20
$ws.Range("C38").Value = 181
$ws.Range("D38").Value = -161
$ws.Range("F38").Value = "We lent more than we borrowed"
$ws.Range("G38").Value = "0.11 : 1"

$ws.Range("B39").Value = 49
$ws.Range("C39").Value = 78
$ws.Range("D39").Value = -29
$ws.Range("F39").Value = "We lent more than we borrowed"
$ws.Range("G39").Value = "0.63 : 1"

$ws.Range("B40").Value = 91
$ws.Range("C40").Value = 131
$ws.Range("D40").Value = -40
$ws.Range("F40").Value = "We lent more than we borrowed"
$ws.Range("G40").Value = "0.69 : 1"

$ws.Range("B41").Value = 5
$ws.Range("C41").Value = 18
$ws.Range("D41").Value = -13
$ws.Range("F41").Value = "We lent more than we borrowed"
$ws.Range("G41").Value = "0.28 : 1"

$ws.Range("B42").Value = 19
$ws.Range("C42").Value = 19
$ws.Range("D42").Value = 0
$ws.Range("G42").Value = "1.00 : 1"

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0

$ws.Range("B44").Value = 67
$ws.Range("C44").Value = 93
$ws.Range("D44").Value = -26
$ws.Range("F44").Value = "We lent more than we borrowed"
$ws.Range("G44").Value = "0.72 : 1"

$ws.Range("B45").Value = 69
$ws.Range("C45").Value = 172
$ws.Range("D45").Value = -103
$ws.Range("F45").Value = "We lent more than we borrowed"
$ws.Range("G45").Value = "0.40 : 1"

$ws.Range("B46").Value = 469
$ws.Range("C46").Value = 626
$ws.Range("D46").Value = -157
$ws.Range("F46").Value = "We lent more than we borrowed"
$ws.Range("G46").Value = "0.75 : 1"

$ws.Range("B47").Value = 892
$ws.Range("C47").Value = 706
$ws.Range("D47").Value = 186
$ws.Range("E47").Value = "We borrowerd more than we lent"
$ws.Range("G47").Value = "1.26 : 1"

$ws.Range("B48").Value = 264
$ws.Range("C48").Value = 630
$ws.Range("D48").Value = -366
$ws.Range("F48").Value = "We lent more than we borrowed"
$ws.Range("G48").Value = "0.42 : 1"

$ws.Range("B49").Value = 422
$ws.Range("C49").Value = 257
$ws.Range("D49").Value = 165
$ws.Range("E49").Value = "We borrowerd more than we lent"
$ws.Range("G49").Value = "1.64 : 1"

$ws.Range("B50").Value = 792
$ws.Range("C50").Value = 540
$ws.Range("D50").Value = 252
$ws.Range("E50").Value = "We borrowerd more than we lent"
$ws.Range("G50").Value = "1.47 : 1"

$ws.Range("B51").Value = 166
$ws.Range("C51").Value = 146
$ws.Range("D51").Value = 20
$ws.Range("E51").Value = "We borrowerd more than we lent"
$ws.Range("G51").Value = "1.14 : 1"

$ws.Range("B52").Value = 436
$ws.Range("C52").Value = 528
$ws.Range("D52").Value = -92
$ws.Range("F52").Value = "We lent more than we borrowed"
$ws.Range("G52").Value = "0.83 : 1"

$ws.Range("B53").Value = 123
$ws.Range("C53").Value = 266
$ws.Range("D53").Value = -143
$ws.Range("F53").Value = "We lent more than we borrowed"
$ws.Range("G53").Value = "0.46 : 1"

$ws.Range("B54").Value = 19
$ws.Range("C54").Value = 231
$ws.Range("D54").Value = -212
$ws.Range("F54").Value = "We lent more than we borrowed"
$ws.Range("G54").Value = "0.08 : 1"

$ws.Range("B55").Value = 321
$ws.Range("C55").Value = 234
$ws.Range("D55").Value = 87
$ws.Range("E55").Value = "We borrowerd more than we lent"
$ws.Range("G55").Value = "1.37 : 1"
